# Fill in the timesheet entry for Thursday 31.5.18 (Wk 12) and let the
# SUBTOTAL formula in the TOTAL row recalc naturally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings must be introduced in the same order as the target
# workbook (Date, Description, Time) so the sharedStrings table lines up.
$ws.Range("A9").Value = "Wk[12] Thursday 31.5.18"
$ws.Range("D9").Value = "Testing"
$ws.Range("B9").Value = "1700 - 2130"
$ws.Range("C9").Value = 4.5

# Move the active selection to D9, matching the last cell touched.
$ws.Range("D9").Select()

$excel.Calculate()
